$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Yoshi Ilands 2 / Radio/Podcast / 8/8/2020 / Hermione Granger / None / New
$ws.Range("A8").Value = "Yoshi Ilands 2"
$ws.Range("B8").Value = "Radio/Podcast"
$ws.Range("C6").Copy($ws.Range("C8"))
$ws.Range("C8").Value = 44051
$ws.Range("D8").Value = "Hermione Granger"
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "New"

# Row 9: Email design for summer dinner / Email Marketing / 7/11/2019 / Hermione Granger / None / New
$ws.Range("A9").Value = "Email design for summer dinner"
$ws.Range("B9").Value = "Email Marketing"
$ws.Range("C7").Copy($ws.Range("C9"))
$ws.Range("C9").Value = 43657
$ws.Range("D9").Value = "Hermione Granger"
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "New"

# Update selection to match the saved workbook view state
$ws.Range("G19").Select()
